$d = $word.ActiveDocument

# Locate the paragraph that currently reads "Baz changes" (stored as two
# runs "Baz chan" / "ges" with a hidden _GoBack bookmark sitting between
# them).
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "Baz changes") {
        $target = $p
    }
}

if ($target -ne $null) {
    $r = $target.Range
    # Range covering just the paragraph's text, excluding the trailing
    # paragraph mark, so the paragraph itself is preserved.
    $textRange = $d.Range($r.Start, $r.End - 1)

    # Replace the paragraph's contents with the new sentence, written as two
    # runs (the sentence, then a separate "." run) followed by the _GoBack
    # bookmark -- which now sits after all of the text instead of in the
    # middle of it. InsertXML (unlike Range.Text/InsertAfter) preserves the
    # run boundaries exactly as supplied instead of coalescing them.
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>I can see how important Version management is in developing software</w:t></w:r><w:r><w:t>.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $textRange.InsertXML($xml)
}

# The document used to end with two consecutive empty paragraphs; drop one
# of them so a single empty paragraph remains before the section break.
$trailing = @()
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "") {
        $trailing += $p
    }
}
if ($trailing.Count -ge 2) {
    $prev = $trailing[$trailing.Count - 2]
    $prev.Range.Delete()
}
